# "some front end changes"
#
# Applies two pure paragraph-mark / bookmark-placement edits (no visible
# text is added or removed anywhere):
#
#   1. Near the top of the document, the lone "_GoBack" bookmark paragraph
#      (third paragraph) loses its bookmark and its paragraph-mark rFonts
#      hint, becoming a bare empty paragraph; the bookmark instead opens
#      the title paragraph ("GameCore开发·"), which also loses its own
#      paragraph-mark rFonts hint.
#   2. Further down, the "如果赢了的赛道...使用" paragraph loses its
#      paragraph-mark rFonts hint, and the empty paragraph right after it
#      (which only carried a pPr/rFonts hint) becomes a bare empty
#      paragraph too.
#
# Paragraph-mark / bookmark placement can't be massaged into place with
# plain Find/Replace, so each affected paragraph is replaced wholesale with
# an equivalent hand-built XML fragment via Range.InsertXML, which parses
# literal OOXML into the target Range (replacing exactly that range, mark
# included - this mirrors real Word's Range.InsertXML behavior).
#
# Paragraphs are located by exact text match (trimmed of the trailing
# paragraph-mark CR) rather than hard-coded indices, so the script keeps
# working even if unrelated paragraphs were inserted/removed elsewhere.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

function Get-ParaIndexByExactText($target) {
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]10)
        if ($t -eq $target) {
            return $i
        }
    }
    return -1
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$emptyXml = '<w:p ' + $wNs + '/>'

# ---------------------------------------------------------------------
# 1. Title paragraph gets the _GoBack bookmark prepended and loses pPr;
#    the paragraph that used to be the bookmark's sole occupant (two
#    paragraphs later: title, sql-reminder, bookmark) becomes empty.
# ---------------------------------------------------------------------
$idxTitle = Get-ParaIndexByExactText("GameCore开发·")
if ($idxTitle -eq -1) {
    throw "Could not locate the 'GameCore开发·' title paragraph"
}

$titleXml = '<w:p ' + $wNs + '>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>GameCore</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>开发</w:t></w:r>' +
    '<w:r w:rsidR="00261F21"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>·</w:t></w:r>' +
    '</w:p>'

$idxBookmarkPara = $idxTitle + 2
$bookmarkParaText = $d.Paragraphs.Item($idxBookmarkPara).Range.Text.TrimEnd([char]13, [char]10)
if ($bookmarkParaText -ne "") {
    throw "Expected an empty paragraph (the old _GoBack holder) at index $idxBookmarkPara"
}

# Replace the title paragraph first...
$d.Paragraphs.Item($idxTitle).Range.InsertXML($titleXml)
# ...then clear out the paragraph that used to hold the bookmark alone.
$d.Paragraphs.Item($idxBookmarkPara).Range.InsertXML($emptyXml)

# ---------------------------------------------------------------------
# 2. The "如果赢了的赛道...使用" paragraph loses its pPr; the empty
#    paragraph right after it also loses its pPr (becomes bare).
# ---------------------------------------------------------------------
$winpoolText = "如果赢了的赛道没有人投注，那就把所有的token计入一个叫做winpool的user名下，等于归入奖池，日后在其他的活动中使用"
$idxWinpool = Get-ParaIndexByExactText($winpoolText)
if ($idxWinpool -eq -1) {
    throw "Could not locate the 'winpool' paragraph"
}

$winpoolXml = '<w:p ' + $wNs + '>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>如果赢了的赛道没有人投注，那就把所有的token计入一个叫做</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>winpool</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>的user名下，等于归入奖池，日后在其他的活动中使用</w:t></w:r>' +
    '</w:p>'

$idxTrailingEmpty = $idxWinpool + 1
$trailingEmptyText = $d.Paragraphs.Item($idxTrailingEmpty).Range.Text.TrimEnd([char]13, [char]10)
if ($trailingEmptyText -ne "") {
    throw "Expected an empty paragraph right after the 'winpool' paragraph at index $idxTrailingEmpty"
}

$d.Paragraphs.Item($idxWinpool).Range.InsertXML($winpoolXml)
$d.Paragraphs.Item($idxTrailingEmpty).Range.InsertXML($emptyXml)
